# change process CDM mapping
#
# Updates the discharge_to_source_value -> discharge_to_concept_id join
# table: the existing mapping rows (2-5) get new source-value labels and
# concept ids, and two additional "Hospital patient_*" categories are
# appended as new rows (6-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Home"
$ws.Range("B2").Value = 8536

$ws.Range("A3").Value = "Patient died"
$ws.Range("B3").Value = 4216643

$ws.Range("A4").Value = "Patient transfer from hospital to hospital"
$ws.Range("B4").Value = 44790567

$ws.Range("A5").Value = "Hospital patient_Ward"
$ws.Range("B5").Value = 4148418

$ws.Range("A6").Value = "Hospital patient_Intensive Care"
$ws.Range("A7").Value = "Hospital patient_Operating room"

# New rows 6 & 7 should carry the same number style as the existing B
# column entries (style index 1 in the original file), so copy the
# format from B5 instead of re-declaring a brand new font/style.
$ws.Range("B5").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B6").Value = 4148418
$ws.Range("B7").Value = 4148418

# Column A needs to widen a bit to fit the longer labels.
$ws.Columns.Item(1).ColumnWidth = 38.4

# Matches the zoomed-in view and post-edit selection seen in the saved file.
$excel.ActiveWindow.Zoom = 115
$ws.Range("A8").Select()

$wb.Save()
